$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "department" column (C) from the single generic college name to the
# specific department/category each course belongs to.
$ws.Range("C2").Value = "Management"
$ws.Range("C3").Value = "Management"
$ws.Range("C4").Value = "Logistics"
$ws.Range("C5").Value = "Logistics"
$ws.Range("C6").Value = "Management"
$ws.Range("C7").Value = "Management"
$ws.Range("C8").Value = "Information Technology"
$ws.Range("C9").Value = "Information Technology"
$ws.Range("C10").Value = "Graduate Studies"
$ws.Range("C11").Value = "Graduate Studies"
$ws.Range("C12").Value = "Packages"
$ws.Range("C13").Value = "Packages"
$ws.Range("C14").Value = "Packages"
$ws.Range("C15").Value = "Packages"
